# "Minor corrections" commit: a few more game entries were appended to the
# tracking sheet, the header/date rows that used to sit at rows 7-8 were
# pushed down by one (gaining the same formatting the rows above already
# use), and the active-cell selection moved along with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteAll / xlPasteValues paste-special codes used below.
$xlPasteAll    = -4104
$xlPasteValues = -4163

# ---------------------------------------------------------------------
# 1) Rows 4 and 5 keep their existing values; they simply pick up the
#    same "ht=15 / customHeight / row style" formatting that row 3 has.
# ---------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 15
$ws.Rows.Item(5).RowHeight = 15

# ---------------------------------------------------------------------
# 2) The old row 7 ("21/21/12/12/12") shifts up to become row 6, and
#    gains an explicit (empty) styled cell in column F. We reuse the
#    already-styled source cells so the text stays text (not numbers)
#    and the existing style index is preserved instead of minting a
#    new one.
# ---------------------------------------------------------------------
$ws.Range("A7:E7").Copy()
$ws.Range("A6:E6").PasteSpecial($xlPasteAll)

$ws.Range("A1").Copy()
$ws.Range("F6").PasteSpecial($xlPasteAll)
$ws.Range("F6").ClearContents()

$ws.Rows.Item(6).RowHeight = 15

# ---------------------------------------------------------------------
# 3) The old row 8 ("дата / дата конец / .../ 1205233811") shifts up to
#    become row 7, now formatted the same way as the surrounding rows
#    (it previously had no explicit style).
# ---------------------------------------------------------------------
$ws.Range("A5:F5").Copy()
$ws.Range("A7:F7").PasteSpecial($xlPasteAll)

$ws.Range("A8:E8").Copy()
$ws.Range("A7:E7").PasteSpecial($xlPasteValues)
$ws.Cells.Item(7, 6).Value = 1205233811

$ws.Rows.Item(7).RowHeight = 15

$ws.Range("A8:F8").ClearContents()

# ---------------------------------------------------------------------
# 4) New rows 8 and 9, formatted like the rows above them.
# ---------------------------------------------------------------------
$ws.Range("A5:F5").Copy()
$ws.Range("A8:F8").PasteSpecial($xlPasteAll)
$ws.Cells.Item(8, 1).Value = "21 января"
$ws.Cells.Item(8, 2).Value = "23 февраля"
$ws.Cells.Item(8, 3).Value = "БЗР"
$ws.Cells.Item(8, 4).Value = "Вадимир"
$ws.Cells.Item(8, 5).Value = "захват точек"
$ws.Cells.Item(8, 6).Value = 1205233811
$ws.Rows.Item(8).RowHeight = 15

$ws.Range("A5:F5").Copy()
$ws.Range("A9:F9").PasteSpecial($xlPasteAll)
$ws.Cells.Item(9, 1).Value = "12 сентября"
$ws.Cells.Item(9, 2).Value = "13 сентября"
$ws.Cells.Item(9, 3).Value = "СТАЛКЕР"
$ws.Cells.Item(9, 4).Value = "Тур"
$ws.Cells.Item(9, 5).Value = "Ролевая"
$ws.Cells.Item(9, 6).Value = "d1m0nnch1k"
$ws.Rows.Item(9).RowHeight = 15

# ---------------------------------------------------------------------
# 5) New rows 10-12: same "s=2" cell styling as above, but without the
#    explicit row height (matches the diff, which leaves these rows
#    without row-level attributes).
# ---------------------------------------------------------------------
$ws.Range("A5:F5").Copy()
$ws.Range("A10:F10").PasteSpecial($xlPasteAll)
$ws.Cells.Item(10, 1).Value = "12 сентября 2222"
$ws.Cells.Item(10, 2).Value = "15 октября 2222"
$ws.Cells.Item(10, 3).Value = "мероприятие"
$ws.Cells.Item(10, 4).Value = "Дроид групп"
$ws.Cells.Item(10, 5).Value = "Удержание"
$ws.Cells.Item(10, 6).Value = "d1m0nnch1k"

$ws.Range("A5:F5").Copy()
$ws.Range("A11:F11").PasteSpecial($xlPasteAll)
$ws.Cells.Item(11, 1).Value = "1 января 2023"
$ws.Cells.Item(11, 2).Value = "3 января 2023"
$ws.Cells.Item(11, 3).Value = "Битва за Мандарины"
$ws.Cells.Item(11, 4).Value = "Вадимир"
$ws.Cells.Item(11, 5).Value = "Охота за конфетами"
$ws.Cells.Item(11, 6).Value = "d1m0nnch1k"

$ws.Range("A5:F5").Copy()
$ws.Range("A12:F12").PasteSpecial($xlPasteAll)
$ws.Cells.Item(12, 1).Value = "1 февраля 2023"
$ws.Cells.Item(12, 2).Value = "4 февраля 2023"
$ws.Cells.Item(12, 3).Value = "Февральский мороз"
$ws.Cells.Item(12, 4).Value = "ПП"
$ws.Cells.Item(12, 5).Value = "встречный бой"
$ws.Cells.Item(12, 6).Value = "d1m0nnch1k"

# ---------------------------------------------------------------------
# 6) New row 13: plain, unstyled cells (no style carried over at all).
# ---------------------------------------------------------------------
$ws.Cells.Item(13, 1).Value = "13 марта 2022"
$ws.Cells.Item(13, 2).Value = "15 марта 2022"
$ws.Cells.Item(13, 3).Value = "название"
$ws.Cells.Item(13, 4).Value = "организатор"
$ws.Cells.Item(13, 5).Value = "игра"
$ws.Cells.Item(13, 6).Value = "d1m0nnch1k"

# ---------------------------------------------------------------------
# 7) Move the active cell/selection the way the author's session ended.
# ---------------------------------------------------------------------
$ws.Range("C15").Select()
